$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the analytical formula in column AV (rows 2-18): the "-2*X*...+4*2*232000*(...)"
# term is replaced by an equivalent but simplified "-4*X*...*(0.08/0.4)/..." term so the
# analytical solution can be validated against the simulation data.
$ws.Range("AV2").Formula = "=(AQ2*X2-4*X2*(1-0.01*P2-2*0.01*AF2)*(0.08/0.4)/(-0.08/0.4*0.01*P2-(2*0.08/0.4+3)*0.01*AF2+0.08/0.4+1))/1000"
$ws.Range("AV3:AV18").Formula = "=(AQ3*X3-4*X3*(1-0.01*P3-2*0.01*AF3)*(0.08/0.4)/(-0.08/0.4*0.01*P3-(2*0.08/0.4+3)*0.01*AF3+0.08/0.4+1))/1000"

# Match the author's final selection: the whole AV column, anchored at AV1.
$ws.Columns("AV").Select()
